# Applies scheduled-runner market data refresh to the Leve profit tables
# (currentAveragePrice* / LevePrice* / LeveProfit* columns) across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3114.2856
$ws.Range("I64").Value = 2850
$ws.Range("J64").Value = 3466.6667
$ws.Range("K64").Value = 2850
$ws.Range("L64").Value = 3466.6667
$ws.Range("M64").Value = -2602
$ws.Range("N64").Value = -3962.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3114.2856
$ws.Range("I67").Value = 2850
$ws.Range("J67").Value = 3466.6667
$ws.Range("K67").Value = 2850
$ws.Range("L67").Value = 3466.6667
$ws.Range("M67").Value = -1992
$ws.Range("N67").Value = -5182.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1140.0741
$ws.Range("I129").Value = 336.4
$ws.Range("J129").Value = 1322.7273
$ws.Range("K129").Value = 1009.2
$ws.Range("L129").Value = 3968.1819
$ws.Range("M129").Value = 3990.8
$ws.Range("N129").Value = -13968.1819

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2377
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7504.963
$ws.Range("I32").Value = 5031.341
$ws.Range("K32").Value = 5031.341
$ws.Range("M32").Value = -4744.341

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 907.7619
$ws.Range("J74").Value = 1693.5834
$ws.Range("L74").Value = 1693.5834
$ws.Range("N74").Value = -3441.5834

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 907.7619
$ws.Range("J77").Value = 1693.5834
$ws.Range("L77").Value = 8467.916999999999
$ws.Range("N77").Value = -17203.917

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1242.0605
$ws.Range("I132").Value = 999.6774
$ws.Range("K132").Value = 2999.0322
$ws.Range("M132").Value = -469.0322000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 308.25
$ws.Range("I22").Value = 80
$ws.Range("J22").Value = 536.5
$ws.Range("K22").Value = 80
$ws.Range("L22").Value = 536.5
$ws.Range("M22").Value = 93
$ws.Range("N22").Value = -882.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 171082.25
$ws.Range("I86").Value = 4817
$ws.Range("K86").Value = 4817
$ws.Range("M86").Value = -3694

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 171082.25
$ws.Range("I89").Value = 4817
$ws.Range("K89").Value = 24085
$ws.Range("M89").Value = -18469

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 988.625
$ws.Range("I94").Value = 421.8
$ws.Range("J94").Value = 1933.3334
$ws.Range("K94").Value = 421.8
$ws.Range("L94").Value = 1933.3334
$ws.Range("M94").Value = 29.19999999999999
$ws.Range("N94").Value = -2835.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4063.8333
$ws.Range("J31").Value = 3032.3333
$ws.Range("L31").Value = 3032.3333
$ws.Range("N31").Value = -3622.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4063.8333
$ws.Range("J34").Value = 3032.3333
$ws.Range("L34").Value = 3032.3333
$ws.Range("N34").Value = -3436.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2560310.2
$ws.Range("I58").Value = 3108054.5
$ws.Range("J58").Value = 4171
$ws.Range("K58").Value = 3108054.5
$ws.Range("L58").Value = 4171
$ws.Range("M58").Value = -3107851.5
$ws.Range("N58").Value = -4577

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1647.3334
$ws.Range("I132").Value = 1088.7778
$ws.Range("K132").Value = 3266.3334
$ws.Range("M132").Value = -736.3334000000004

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2560310.2
$ws.Range("I136").Value = 3108054.5
$ws.Range("J136").Value = 4171
$ws.Range("K136").Value = 9324163.5
$ws.Range("L136").Value = 12513
$ws.Range("M136").Value = -9321613.5
$ws.Range("N136").Value = -17613

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 380.66666
$ws.Range("I2").Value = 361.66666
$ws.Range("J2").Value = 399.66666
$ws.Range("K2").Value = 2169.99996
$ws.Range("L2").Value = 2397.99996
$ws.Range("M2").Value = -2056.99996
$ws.Range("N2").Value = -2623.99996

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 115.5
$ws.Range("I12").Value = 70.8
$ws.Range("J12").Value = 132.6923
$ws.Range("K12").Value = 212.4
$ws.Range("L12").Value = 398.0769
$ws.Range("M12").Value = -39.39999999999998
$ws.Range("N12").Value = -744.0769

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 37232.25
$ws.Range("I129").Value = 730
$ws.Range("J129").Value = 43673.824
$ws.Range("K129").Value = 2190
$ws.Range("L129").Value = 131021.472
$ws.Range("M129").Value = 2810
$ws.Range("N129").Value = -141021.472

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4703.8
$ws.Range("I70").Value = 4250
$ws.Range("J70").Value = 5006.3335
$ws.Range("K70").Value = 4250
$ws.Range("L70").Value = 5006.3335
$ws.Range("M70").Value = -3980
$ws.Range("N70").Value = -5546.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4703.8
$ws.Range("I73").Value = 4250
$ws.Range("J73").Value = 5006.3335
$ws.Range("K73").Value = 4250
$ws.Range("L73").Value = 5006.3335
$ws.Range("M73").Value = -3314
$ws.Range("N73").Value = -6878.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2229.6667
$ws.Range("I80").Value = 2345
$ws.Range("K80").Value = 2345
$ws.Range("M80").Value = -1347

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2229.6667
$ws.Range("I83").Value = 2345
$ws.Range("K83").Value = 11725
$ws.Range("M83").Value = -6733

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2406908
$ws.Range("I132").Value = 2960894.5
$ws.Range("J132").Value = 6299.3335
$ws.Range("K132").Value = 8882683.5
$ws.Range("L132").Value = 18898.0005
$ws.Range("M132").Value = -8880153.5
$ws.Range("N132").Value = -23958.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 12532.333
$ws.Range("I40").Value = 12489.2
$ws.Range("K40").Value = 12489.2
$ws.Range("M40").Value = -12353.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 12422.223
$ws.Range("I122").Value = 13816.667
$ws.Range("K122").Value = 41450.001
$ws.Range("M122").Value = -39000.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2268.7754
$ws.Range("I132").Value = 1388.6364
$ws.Range("K132").Value = 4165.9092
$ws.Range("M132").Value = -1635.9092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 557.3200000000001
$ws.Range("I113").Value = 409.3125
$ws.Range("J113").Value = 820.44446
$ws.Range("K113").Value = 1227.9375
$ws.Range("L113").Value = 2461.33338
$ws.Range("M113").Value = 942.0625
$ws.Range("N113").Value = -6801.33338

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1503.1613
$ws.Range("I132").Value = 1129.6786
$ws.Range("K132").Value = 3389.0358
$ws.Range("M132").Value = -859.0357999999997

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 24156374
$ws.Range("I136").Value = 50506480
$ws.Range("J136").Value = 2110.4167
$ws.Range("K136").Value = 151519440
$ws.Range("L136").Value = 6331.250100000001
$ws.Range("M136").Value = -151516890
$ws.Range("N136").Value = -11431.2501
